$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Helper: write literal text into a cell without Excel's automatic
# "looks like a date" conversion, and without disturbing the
# destination cell's existing style/number format.
# Strategy: stage the text as Text-formatted in an out-of-the-way
# scratch cell, copy it, then paste only the VALUE into the target
# cell (PasteSpecial xlPasteValues = -4163 copies the value verbatim
# as text without touching the destination's formatting).
$scratch = $ws.Range("Z100")

function Set-LiteralText($addr, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# Row 3
$ws.Range("H3").Value = 35
Set-LiteralText "I3" "04-Nov-2025"

# Row 4
$ws.Range("H4").Value = -151
Set-LiteralText "I4" "04-Nov-2025"

# Row 5
$ws.Range("H5").Value = -100
Set-LiteralText "I5" "04-Nov-2025"

# Row 6
$ws.Range("H6").Value = 300
Set-LiteralText "I6" "04-Nov-2025"

# Clean up the scratch cell so it leaves no trace in the sheet.
$scratch.Clear()
